$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.593.98'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.08%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.517.88'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.22%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '202.61'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.49%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '555.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.96%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.506.35'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.37%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.609'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.77%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.657'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.14%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '63.52'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +12.59%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.144'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -6.30%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000272'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.84%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.90'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.076.59'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.38%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.512.45'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.55%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.81%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.40%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.331.61'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.46%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.85'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -5.71%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.28%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '393.41'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.32%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.19'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -7.26%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.01'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.78%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.28'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.12%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.90'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.05%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.83'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.23%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.23'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.18%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.87'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.41%  '

# Row 30
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '716.18'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.96%  '

# Row 31
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.09'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.10%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.10'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -13.35%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.77'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.88%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.94'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.18%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.112'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.70%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.52'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -10.01%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.08%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.398'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.97%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.132'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.63%  '

# Row 40
$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.00'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.32%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.16%  '

# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.056.03'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.90%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0686'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -12.83%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.60'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -9.93%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.77'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.55%  '

# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.72'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -9.39%  '

# Row 47
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0406'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.45%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.128'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.25%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.79'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.22%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.27'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.92%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.89'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.86%  '
